$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 7.613107666666667
$ws.Range("H2").Value = 22.839323
$ws.Range("I2").Value = 0.08102996839946881
$ws.Range("J2").Value = 0.0810299683994688
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 1279.83821531341
$ws.Range("R2").Value = 11518.54393782069
$ws.Range("S2").Value = 0.02418081940666305
$ws.Range("T2").Value = 0.02418081940666305
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 7.613107666666667
$ws.Range("H3").Value = 22.839323
$ws.Range("I3").Value = 0.08102996839946881
$ws.Range("J3").Value = 0.0810299683994688
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("Q3").Value = 1240.984022468373
$ws.Range("R3").Value = 11168.85620221536
$ws.Range("S3").Value = 0.02344672176124511
$ws.Range("T3").Value = 0.02344672176124511
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 7.613107666666667
$ws.Range("H4").Value = 22.839323
$ws.Range("I4").Value = 0.08102996839946881
$ws.Range("J4").Value = 0.0810299683994688
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 1263.726615860063
$ws.Range("R4").Value = 11373.53954274057
$ws.Range("S4").Value = 0.02387641243391264
$ws.Range("T4").Value = 0.02387641243391263
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 7.613107666666667
$ws.Range("H5").Value = 22.839323
$ws.Range("I5").Value = 0.08102996839946881
$ws.Range("J5").Value = 0.0810299683994688
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 504.1912588913966
$ws.Range("R5").Value = 4537.721330022569
$ws.Range("S5").Value = 0.009526014797648011
$ws.Range("T5").Value = 0.009526014797648009
$ws.Range("I6").Value = 0.7831116101658118
$ws.Range("J6").Value = 0.7831116101658117
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.2984182258032519
$ws.Range("P6").Value = 0.298418225803252
$ws.Range("Q6").Value = 12368.95663842309
$ws.Range("R6").Value = 111320.6097458078
$ws.Range("S6").Value = 0.2336947773116094
$ws.Range("T6").Value = 0.2336947773116094
$ws.Range("I7").Value = 0.7831116101658118
$ws.Range("J7").Value = 0.7831116101658117
$ws.Range("O7").Value = 0.2893586437755394
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("S7").Value = 0.2266001134424582
$ws.Range("T7").Value = 0.2266001134424582
$ws.Range("I8").Value = 0.7831116101658118
$ws.Range("J8").Value = 0.7831116101658117
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 12213.24658645744
$ws.Range("R8").Value = 109919.2192781169
$ws.Range("S8").Value = 0.2307528455882615
$ws.Range("T8").Value = 0.2307528455882615
$ws.Range("I9").Value = 0.7831116101658118
$ws.Range("J9").Value = 0.7831116101658117
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.1175616254801657
$ws.Range("P9").Value = 0.1175616254801657
$ws.Range("Q9").Value = 4872.740744948354
$ws.Range("R9").Value = 43854.66670453519
$ws.Range("S9").Value = 0.09206387382348265
$ws.Range("T9").Value = 0.09206387382348265
$ws.Range("G10").Value = 12.72068066666667
$ws.Range("H10").Value = 38.162042
$ws.Range("I10").Value = 0.1353923256534006
$ws.Range("J10").Value = 0.1353923256534005
$ws.Range("M10").Value = 168.1098273333333
$ws.Range("N10").Value = 504.329482
$ws.Range("O10").Value = 0.2984182258032519
$ws.Range("P10").Value = 0.298418225803252
$ws.Range("Q10").Value = 2138.471430435805
$ws.Range("R10").Value = 19246.24287392224
$ws.Range("S10").Value = 0.04040353760886391
$ws.Range("T10").Value = 0.04040353760886391
$ws.Range("G11").Value = 12.72068066666667
$ws.Range("H11").Value = 38.162042
$ws.Range("I11").Value = 0.1353923256534006
$ws.Range("J11").Value = 0.1353923256534005
$ws.Range("O11").Value = 0.2893586437755394
$ws.Range("P11").Value = 0.2893586437755394
$ws.Range("Q11").Value = 2073.550270591077
$ws.Range("R11").Value = 18661.95243531969
$ws.Range("S11").Value = 0.03917693972868416
$ws.Range("T11").Value = 0.03917693972868416
$ws.Range("G12").Value = 12.72068066666667
$ws.Range("H12").Value = 38.162042
$ws.Range("I12").Value = 0.1353923256534006
$ws.Range("J12").Value = 0.1353923256534005
$ws.Range("M12").Value = 165.99353
$ws.Range("N12").Value = 497.98059
$ws.Range("O12").Value = 0.294661504941043
$ws.Range("P12").Value = 0.294661504941043
$ws.Range("Q12").Value = 2111.550687862753
$ws.Range("R12").Value = 19003.95619076478
$ws.Range("S12").Value = 0.03989490643449879
$ws.Range("T12").Value = 0.03989490643449878
$ws.Range("G13").Value = 12.72068066666667
$ws.Range("H13").Value = 38.162042
$ws.Range("I13").Value = 0.1353923256534006
$ws.Range("J13").Value = 0.1353923256534005
$ws.Range("M13").Value = 66.22673433333334
$ws.Range("N13").Value = 198.680203
$ws.Range("O13").Value = 0.1175616254801657
$ws.Range("P13").Value = 0.1175616254801657
$ws.Range("Q13").Value = 842.4491390505029
$ws.Range("R13").Value = 7582.042251454526
$ws.Range("S13").Value = 0.0159169418813537
$ws.Range("T13").Value = 0.0159169418813537
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.04379166666666667
$ws.Range("H14").Value = 0.131375
$ws.Range("I14").Value = 0.0004660957813189216
$ws.Range("J14").Value = 0.0004660957813189215
$ws.Range("M14").Value = 168.1098273333333
$ws.Range("N14").Value = 504.329482
$ws.Range("O14").Value = 0.2984182258032519
$ws.Range("P14").Value = 0.298418225803252
$ws.Range("Q14").Value = 7.361809521972223
$ws.Range("R14").Value = 66.25628569775
$ws.Range("S14").Value = 0.0001390914761155731
$ws.Range("T14").Value = 0.0001390914761155731
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.04379166666666667
$ws.Range("H15").Value = 0.131375
$ws.Range("I15").Value = 0.0004660957813189216
$ws.Range("J15").Value = 0.0004660957813189215
$ws.Range("O15").Value = 0.2893586437755394
$ws.Range("P15").Value = 0.2893586437755394
$ws.Range("Q15").Value = 7.138314736902779
$ws.Range("R15").Value = 64.24483263212501
$ws.Range("S15").Value = 0.0001348688431519436
$ws.Range("T15").Value = 0.0001348688431519435
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.04379166666666667
$ws.Range("H16").Value = 0.131375
$ws.Range("I16").Value = 0.0004660957813189216
$ws.Range("J16").Value = 0.0004660957813189215
$ws.Range("M16").Value = 165.99353
$ws.Range("N16").Value = 497.98059
$ws.Range("O16").Value = 0.294661504941043
$ws.Range("P16").Value = 0.294661504941043
$ws.Range("Q16").Value = 7.269133334583334
$ws.Range("R16").Value = 65.42220001125001
$ws.Range("S16").Value = 0.0001373404843701047
$ws.Range("T16").Value = 0.0001373404843701047
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.04379166666666667
$ws.Range("H17").Value = 0.131375
$ws.Range("I17").Value = 0.0004660957813189216
$ws.Range("J17").Value = 0.0004660957813189215
$ws.Range("M17").Value = 66.22673433333334
$ws.Range("N17").Value = 198.680203
$ws.Range("O17").Value = 0.1175616254801657
$ws.Range("P17").Value = 0.1175616254801657
$ws.Range("Q17").Value = 2.900179074347223
$ws.Range("R17").Value = 26.10161166912501
$ws.Range("S17").Value = 0.00005479497768130025
$ws.Range("T17").Value = 0.00005479497768130025
